$wb = $excel.ActiveWorkbook

$settings = $wb.Worksheets.Item("settings")
$settings.Range("C1").Value = "display.title.text"
[void]$settings.Range("C2").Select()

$survey = $wb.Worksheets.Item("survey")
$survey.Range("C1").Value = "display.prompt.text"
$survey.Range("D1").Value = "display.hint.text"
[void]$survey.Range("D2").Select()
$survey.Activate()
